$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "Dee165"
$ws.Range("B2").Value = 23071031
$ws.Range("C2").Value = "jadeja94"
$ws.Range("D2").Value = "tY73%&Sc"
$ws.Range("E2").Value = "MS"
$ws.Range("F2").Value = "Jadeja"
$ws.Range("G2").Value = "Bhaiu"
